$d = $word.ActiveDocument
$rng = $d.Content
$newBodyXml = @'
<w:p w:rsidP="009168BC" w:rsidR="00DE5A1F" w:rsidRDefault="00DE5A1F"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r w:rsidRPr="00DE5A1F"><w:t>This template demonstrates the use of bookmarks. It creates two links to the same bookmark.</w:t></w:r></w:p><w:p w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>dangling reference for bookmark bookmark1</w:t></w:r><w:r/><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:instrText xml:space="preserve"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p w:rsidP="009168BC" w:rsidR="00C52979" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t>Test</w:t></w:r><w:r w:rsidR="00C52979"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>bookmark</w:t></w:r><w:r w:rsidR="00C52979"><w:t xml:space="preserve"> : </w:t></w:r><w:r w:rsidR="00C52979"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="002D1E44"><w:instrText>m</w:instrText></w:r><w:r w:rsidR="002848A5"><w:instrText>:</w:instrText></w:r><w:r><w:instrText>bookmark</w:instrText></w:r><w:r w:rsidR="003D27D6"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00B71CCB"><w:instrText>self.</w:instrText></w:r><w:r w:rsidR="00C52979"><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>Invalid bookmark statement: Expression "self." is invalid: missing feature access or service call</w:t></w:r><w:r w:rsidR="003D27D6"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="003D27D6"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="002D1E44"><w:instrText>m</w:instrText></w:r><w:r w:rsidR="003D27D6"><w:instrText>:end</w:instrText></w:r><w:r><w:instrText>bookmark</w:instrText></w:r><w:r w:rsidR="003D27D6"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="003D27D6"><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p w:rsidP="00E02A2B" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Test link after bookmark : </w:t></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>dangling reference for bookmark bookmark1</w:t></w:r><w:r/><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:instrText xml:space="preserve"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="36D3B6B7C354B3F26DAC558E4E600907"><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00D0546C"><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w:rsidP="00DE5A1F" w:rsidR="007A2DC4" w:rsidRDefault="00DE5A1F"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00DE5A1F"><w:t>End of demonstration.</w:t></w:r></w:p>
'@
$rng.InsertXML($newBodyXml)
Write-Output "rsid update applied"
